# Raven9A_test_1 bug-fix: corrected recorded answers on "question_answers"
# and corresponding computed outputs on "outputs".

$wb = $excel.ActiveWorkbook

# ---- question_answers: column B holds the answers as text values ----
$ws1 = $wb.Worksheets.Item("question_answers")

# keep these as text cells (matches original inlineStr/"Text" cell type);
# split into the contiguous runs of rows that actually change so the
# untouched rows (11, 28, 36) keep their original formatting
$ws1.Range("B2:B10").NumberFormat = "@"
$ws1.Range("B12:B27").NumberFormat = "@"
$ws1.Range("B29:B35").NumberFormat = "@"
$ws1.Range("B37:B49").NumberFormat = "@"

$ws1.Range("B2").Value = "4"
$ws1.Range("B3").Value = "7"
$ws1.Range("B4").Value = "1"
$ws1.Range("B5").Value = "4"
$ws1.Range("B6").Value = "5"
$ws1.Range("B7").Value = "5"
$ws1.Range("B8").Value = "5"
$ws1.Range("B9").Value = "4"
$ws1.Range("B10").Value = "8"
$ws1.Range("B12").Value = "6"
$ws1.Range("B13").Value = "7"
$ws1.Range("B14").Value = "1"
$ws1.Range("B15").Value = "1"
$ws1.Range("B16").Value = "7"
$ws1.Range("B17").Value = "5"
$ws1.Range("B18").Value = "1"
$ws1.Range("B19").Value = "8"
$ws1.Range("B20").Value = "2"
$ws1.Range("B21").Value = "3"
$ws1.Range("B22").Value = "5"
$ws1.Range("B23").Value = "5"
$ws1.Range("B24").Value = "8"
$ws1.Range("B25").Value = "6"
$ws1.Range("B26").Value = "6"
$ws1.Range("B27").Value = "7"
$ws1.Range("B29").Value = "3"
$ws1.Range("B30").Value = "7"
$ws1.Range("B31").Value = "5"
$ws1.Range("B32").Value = "5"
$ws1.Range("B33").Value = "2"
$ws1.Range("B34").Value = "1"
$ws1.Range("B35").Value = "1"
$ws1.Range("B37").Value = "5"
$ws1.Range("B38").Value = "4"
$ws1.Range("B39").Value = "3"
$ws1.Range("B40").Value = "2"
$ws1.Range("B41").Value = "1"
$ws1.Range("B42").Value = "1"
$ws1.Range("B43").Value = "5"
$ws1.Range("B44").Value = "2"
$ws1.Range("B45").Value = "6"
$ws1.Range("B46").Value = "1"
$ws1.Range("B47").Value = "2"
$ws1.Range("B48").Value = "6"
$ws1.Range("B49").Value = "2"

# ---- outputs: recomputed pre/raw scores, and report/level collapsed ----
$ws2 = $wb.Worksheets.Item("outputs")

$ws2.Range("B2").Value = 7
$ws2.Range("B3").Value = 11

# "report" row (B6) no longer holds descriptive text -- it becomes 0,
# and the old "level" row (row 7) is removed entirely.
$ws2.Rows.Item(7).Delete()
$ws2.Range("B6").Value = 0
